$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.401.10"
$ws.Range("E2").Value = "  +1.44%  "
$ws.Range("D3").Value = "1.692.41"
$ws.Range("E4").Value = "  +0.75%  "
$ws.Range("E5").Value = "  +1.48%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.5553"
$ws.Range("E6").Value = "  +8.73%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.010"
$ws.Range("E7").Value = "  +0.73%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2721"
$ws.Range("E8").Value = "  +1.55%  "
$ws.Range("E9").Value = "  +1.77%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "22.17"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07668"
$ws.Range("E11").Value = "  +3.13%  "
$ws.Range("B12").Value = "Polkadot"
$ws.Range("C12").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "4.564"
$ws.Range("E12").Value = "  +1.28%  "
$ws.Range("B13").Value = "Polygon"
$ws.Range("C13").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.5829"
$ws.Range("E13").Value = "  +0.47%  "
$ws.Range("B14").Value = "WrappedEther"
$ws.Range("C14").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D14").Value = "1.553.03"
$ws.Range("E14").Value = "  -6.85%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.000008467"
$ws.Range("E15").Value = "  -0.08%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "65.34"
$ws.Range("E16").Value = "  +2.07%  "
$ws.Range("D17").Value = "26.502.68"
$ws.Range("E17").Value = "  +2.59%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "4.969"
$ws.Range("E18").Value = "  +0.98%  "
$ws.Range("E19").Value = "  +0.70%  "
$ws.Range("E20").Value = "  +1.88%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "190.43"
$ws.Range("E21").Value = "  +0.67%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.255"
$ws.Range("E22").Value = "  +1.11%  "
$ws.Range("E23").Value = "  +0.67%  "
$ws.Range("E24").Value = "  +3.46%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.1308"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "7.914"
$ws.Range("E26").Value = "  +4.14%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "15.77"
$ws.Range("E27").Value = "  +0.77%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.422"
$ws.Range("E28").Value = "  +7.02%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.06333"
$ws.Range("E29").Value = "  -4.39%  "
$ws.Range("E30").Value = "  +1.53%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "3.598"
$ws.Range("E31").Value = "  +1.25%  "
$ws.Range("E32").Value = "  +2.54%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.678"
$ws.Range("E33").Value = "  +1.02%  "
$ws.Range("E34").Value = "  +2.75%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.6219"
$ws.Range("E35").Value = "  +0.65%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.406"
$ws.Range("E36").Value = "  +1.62%  "
$ws.Range("E37").Value = "  +1.51%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "6.235"
$ws.Range("E38").Value = "  -0.94%  "
$ws.Range("D39").Value = "1.125.32"
$ws.Range("E39").Value = "  +2.68%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.01643"
$ws.Range("E40").Value = "  +3.20%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.8817"
$ws.Range("E41").Value = "  +1.80%  "
$ws.Range("E42").Value = "  +0.81%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "100.89"
$ws.Range("E43").Value = "  -0.62%  "
$ws.Range("D44").Value = "1.843.99"
$ws.Range("E44").Value = "  +1.63%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.00000000110"
$ws.Range("E45").Value = "  -4.38%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "57.58"
$ws.Range("E46").Value = "  +2.45%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "8.216"
$ws.Range("E47").Value = "  +1.48%  "
$ws.Range("E48").Value = "  +0.15%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.05282"
$ws.Range("E49").Value = "  +1.06%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.4302"
$ws.Range("E50").Value = "  +0.46%  "
$ws.Range("E51").Value = "  +1.70%  "
